$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.923.08"
$ws.Range("E2").Value = "  -4.95%  "
$ws.Range("D3").Value = "3.282.28"
$ws.Range("E3").Value = "  -5.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "563.69"
$ws.Range("E5").Value = "  -3.58%  "
$ws.Range("E6").Value = "  -3.36%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "3.283.10"
$ws.Range("E8").Value = "  -5.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.473"
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.27"
$ws.Range("E10").Value = "  -4.58%  "
$ws.Range("E11").Value = "  -4.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.370"
$ws.Range("E12").Value = "  -4.26%  "
$ws.Range("D13").Value = "3.847.16"
$ws.Range("E13").Value = "  -5.41%  "
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").Value = "3.286.55"
$ws.Range("E15").Value = "  -5.42%  "
$ws.Range("E16").Value = "  -5.80%  "
$ws.Range("D17").Value = "61.032.31"
$ws.Range("E17").Value = "  -4.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "24.19"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("E19").Value = "  -1.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.55"
$ws.Range("E20").Value = "  -2.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.90"
$ws.Range("E21").Value = "  -10.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "348.88"
$ws.Range("E22").Value = "  -9.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.550"
$ws.Range("E23").Value = "  -3.40%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "3.417.28"
$ws.Range("E25").Value = "  -5.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "69.18"
$ws.Range("E26").Value = "  -7.38%  "
$ws.Range("E27").Value = "  -4.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.12"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  -1.37%  "
$ws.Range("E31").Value = "  -2.15%  "
$ws.Range("E32").Value = "  -6.07%  "
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("E34").Value = "  -2.36%  "
$ws.Range("D35").Value = "3.315.57"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.52"
$ws.Range("E36").Value = "  -1.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.19"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.77"
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "160.12"
$ws.Range("E39").Value = "  -1.49%  "
$ws.Range("E40").Value = "  -2.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0748"
$ws.Range("E41").Value = "  -3.44%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.16"
$ws.Range("E43").Value = "  -0.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.29"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.737"
$ws.Range("E45").Value = "  -7.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.13"
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("E47").Value = "  -4.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.20"
$ws.Range("E48").Value = "  -5.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.65"
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.855"
$ws.Range("E50").Value = "  -5.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.27"
$ws.Range("E51").Value = "  +4.59%  "
